# Commit: exclusão do slide que menciona o jogo da WHIP como exemplo de aplicação
# -> Delete the slide that introduces the "Jogo da WHIP" example (title "Jogo Exemplo").

$p = $ppt.ActivePresentation

# Locate the slide whose title is "Jogo Exemplo" (the slide introducing the
# "Jogo da WHIP" example) instead of assuming a fixed index, so the script is
# resilient to any prior reordering.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Jogo Exemplo") {
                $targetIndex = $i
                break
            }
        }
    }
    if ($targetIndex -ne -1) {
        break
    }
}

if ($targetIndex -eq -1) {
    # Fallback: the slide is originally the 9th slide in the deck.
    $targetIndex = 9
}

$p.Slides.Item($targetIndex).Delete()
